$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new nombre_aides (col C), new montant_total (col E)
$changes = @(
    @{Row=57;  C=3711;   E=138344178},
    @{Row=91;  C=151109; E=482222453},
    @{Row=92;  C=409037; E=1593885439},
    @{Row=93;  C=209517; E=1308009836},
    @{Row=94;  C=94165;  E=916632725},
    @{Row=95;  C=50733;  E=930961320},
    @{Row=96;  C=17249;  E=789741036},
    @{Row=104; C=135227; E=272141544},
    @{Row=105; C=8171;   E=16876204},
    @{Row=106; C=18338;  E=41287345},
    @{Row=141; C=80472;  E=280708050},
    @{Row=184; C=68735;  E=134180601}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 3).Value = $change.C
    $ws.Cells.Item($change.Row, 5).Value = $change.E
}

$wb.Save()
